$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PIR sheet: append rows 157-170 (Bathroom / No Motion / Inactive)
# ---------------------------------------------------------------------------
$wsPir = $wb.Worksheets.Item("PIR")

$pirTimes = @(
    "17:09:05","17:09:06","17:09:07","17:09:12","17:09:17","17:09:22",
    "17:09:27","17:09:32","17:09:37","17:09:42","17:09:47","17:09:52",
    "17:09:57","17:10:02"
)

$pirStartRow = 157
$pirEndRow = $pirStartRow + $pirTimes.Length - 1
$wsPir.Range("A$($pirStartRow):A$($pirEndRow)").NumberFormat = "@"

$r = $pirStartRow
foreach ($t in $pirTimes) {
    $wsPir.Cells.Item($r, 1).Value = "2026-01-30"
    $wsPir.Cells.Item($r, 2).Value = $t
    $wsPir.Cells.Item($r, 3).Value = "17:00"
    $wsPir.Cells.Item($r, 4).Value = "Bathroom"
    $wsPir.Cells.Item($r, 5).Value = "No Motion"
    $wsPir.Cells.Item($r, 6).Value = "Inactive"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Humidity sheet: append rows 102-114 (Bathroom / humidity% / Active)
# ---------------------------------------------------------------------------
$wsHum = $wb.Worksheets.Item("Humidity")

$humRows = @(
    @("17:09:05","87.4%"),
    @("17:09:06","87.4%"),
    @("17:09:07","87.4%"),
    @("17:09:12","87.4%"),
    @("17:09:17","86.5%"),
    @("17:09:22","87.4%"),
    @("17:09:28","87.4%"),
    @("17:09:38","86.5%"),
    @("17:09:43","87.4%"),
    @("17:09:48","87.4%"),
    @("17:09:53","87.4%"),
    @("17:09:58","86.5%"),
    @("17:10:03","87.4%")
)

$humStartRow = 102
$humEndRow = $humStartRow + $humRows.Length - 1
$wsHum.Range("A$($humStartRow):A$($humEndRow)").NumberFormat = "@"
$wsHum.Range("E$($humStartRow):E$($humEndRow)").NumberFormat = "@"

$r = $humStartRow
foreach ($row in $humRows) {
    $wsHum.Cells.Item($r, 1).Value = "2026-01-30"
    $wsHum.Cells.Item($r, 2).Value = $row[0]
    $wsHum.Cells.Item($r, 3).Value = "17:00"
    $wsHum.Cells.Item($r, 4).Value = "Bathroom"
    $wsHum.Cells.Item($r, 5).Value = $row[1]
    $wsHum.Cells.Item($r, 6).Value = "Active"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# mmWave sheet: append row 41 (Living Room / FALL_DETECTED / EMERGENCY)
# ---------------------------------------------------------------------------
$wsMmWave = $wb.Worksheets.Item("mmWave")

$wsMmWave.Range("A41").NumberFormat = "@"
$wsMmWave.Cells.Item(41, 1).Value = "2026-01-30"
$wsMmWave.Cells.Item(41, 2).Value = "17:09:06"
$wsMmWave.Cells.Item(41, 3).Value = "17:00"
$wsMmWave.Cells.Item(41, 4).Value = "Living Room"
$wsMmWave.Cells.Item(41, 5).Value = "FALL_DETECTED"
$wsMmWave.Cells.Item(41, 6).Value = "EMERGENCY"

Write-Output "Applied PIR($pirStartRow-$pirEndRow), Humidity($humStartRow-$humEndRow), mmWave(41) updates"
